# Replace the stale "mean.default" warning / "[1] NA" SourceCode block with
# the lm() summary output, as a single merged paragraph (matching the
# target OOXML structure exactly via a raw XML insert), then remove the
# now-redundant trailing paragraph that used to hold "## [1] NA".

$d = $word.ActiveDocument

# Locate the paragraph that currently contains the "mean.default" warning.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("mean.default")) {
        $targetIndex = $i
        break
    }
}

$p1 = $d.Paragraphs.Item($targetIndex)
$r = $d.Range($p1.Range.Start, $p1.Range.End)

$xmlFragment = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## Call:</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## lm(formula = MEDV ~ CRIME + ZL + NR_PROP + CHR_V + NOX + ROOM + </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##     AGE + DIS + HWY + TAX + PT_RATIO + B + L_PER, data = housing)</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## Coefficients:</w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## (Intercept)        CRIME           ZL      NR_PROP        CHR_V  </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##   3.646e+01   -1.080e-01    4.642e-02    2.056e-02    2.687e+00  </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##         NOX         ROOM          AGE          DIS          HWY  </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##  -1.777e+01    3.810e+00    6.922e-04   -1.476e+00    3.060e-01  </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##         TAX     PT_RATIO            B        L_PER  </w:t></w:r><w:r><w:br w:type="textWrapping"/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##  -1.233e-02   -9.527e-01    9.312e-03   -5.248e-01</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xmlFragment)

# The following paragraph (originally "## [1] NA") is no longer needed;
# its content has been folded into the paragraph above. Delete it outright.
$p1 = $d.Paragraphs.Item($targetIndex)
$p2 = $d.Paragraphs.Item($targetIndex + 1)
if ($p2.Range.Text.Contains("[1] NA")) {
    $p2.Range.Delete()
}
